$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeployNginx")

# --- Fill in the new "Deploy with subdirection" row (row 11) ---
$ws.Range("A11").Value = "Deploy with subdirection"

$nginxConf = @'
server{
    listen  80;
    server_name test.anhduong.us;
        location /icon {
        rewrite ^/icon/?(.*) /$1 break;
        proxy_pass         http://localhost:5002;
        proxy_http_version 1.1;
        proxy_set_header   Upgrade $http_upgrade;
        proxy_set_header   Connection keep-alive;
        proxy_set_header   Host $host;
        proxy_cache_bypass $http_upgrade;
        proxy_set_header   X-Forwarded-For $proxy_add_x_forwarded_for;
        proxy_set_header   X-Forwarded-Proto $scheme;
    }
}

'@
$ws.Range("C11").Value = $nginxConf

# Match the formatting already used for the other descriptive rows in this
# table (Arial 14, wrap text, left/center aligned) instead of the plain
# default style the blank row started with.
foreach ($addr in @("A11", "C11")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 14
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
}

# Row grows tall to fit the wrapped, multi-line nginx config text.
$ws.Rows.Item(11).RowHeight = 297.5

# --- Make "DeployNginx" the active sheet / tab, with B11 selected and ---
# --- the view scrolled so row 9 is at the top (matches the authored view). ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("B11").Select() | Out-Null
